$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cell D1 "Ano"
$ws.Range("D1").Value = "Ano"

# Match the formatting of the other header cells (A1:C1) by copying their format
$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats

# Fill D2:D10 with the reference year range "2023/2012"
for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 4).Value = "2023/2012"
}
